# This script reproduces the commit: "refactor a lot and changed the
# folders structure" applied to the IoT dataset sheet.
#
# Functional changes made to the worksheet data:
#   1. Column A ("time") values are rewritten - the dataset used to be one
#      continuously increasing time series (0..17700 in steps of 300); it
#      is now split into several shorter experiment runs, each of which
#      restarts its own relative clock at 1200 and increases in steps of
#      300 seconds.
#   2. The header cell A1 ("time") loses the bold/bordered header style
#      that is still applied to B1:L1 (it becomes a plain/default cell).
#   3. The active selection on the sheet changes from G7 to N5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for A2:A61 ("time" column), reflecting the new, shorter,
# repeating experiment runs.
$timeValues = @(
    1200, 1500, 1800, 2100, 2400, 2700, 3000, 3300, 3600, 3900, 4200, 4500,
    1200, 1500, 1800, 2100, 2400, 2700, 3000, 3300, 3600, 3900, 4200, 4500,
    1200, 1500, 1800, 2100, 2400, 2700, 3000, 3300, 3600, 3900, 4200, 4500,
    4800, 5100, 5400, 5700, 6000, 6300, 6600, 6900, 7200, 7500, 7800, 8100,
    1200, 1500, 1800, 2100, 2400, 2700, 3000, 3300, 3600, 3900, 4200, 4500
)

$startRow = 2
for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $timeValues[$i]
}

# The "time" header (A1) loses the bold/bordered style that the rest of
# the header row (B1:L1) still has.
$ws.Range("A1").Style = "Normal"

# Update the active cell/selection saved with the sheet view.
$ws.Range("N5").Select()
